$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New measurement data for rows 2-5 (columns A..AH), replacing prior contents.
$data = @(
  @(45166.50694444445,15.374,10.112,3.717,32.732,25.126,11.943,36.363,18.615,7.558,11.183,12.932,13.607,3.86,12.031,16.608,10.57,3.118,1.764,175.75,33.355,11.105,21.604,11.537,2.926,18.421,9.808999999999999,8.874000000000001,10.611,13.63,3.312,32.661,5.894,13.883),
  @(45166.51388888889,0.216,0.626,1.24,0.5590000000000001,0.006,0,7.509,0.569,0,0.127,0.127,0,0,0.376,0.515,0.664,1.3,0.498,0,1.688,0.347,0.878,0.511,0.6909999999999999,3.197,0.226,0.196,0.246,0.286,1.115,7.805,0.135,0.378),
  @(45166.52083333334,7.497,5.964,1.016,16.454,13.174,6.205,20.396,9.269,3.887,6.007,6.528,6.812,1.834,6.015,8.301,5.271,0.921,0.527,84.62,16.608,5.552,10.808,5.903,1.217,9.776,4.832,4.256,5.017,6.983,0.711,18.36,3.144,6.866),
  @(45166.52777777778,19.07,14.58,1.19,41.63,34.03,15.32,54.86,23.23,10.26,15.41,16.65,17.54,4.77,15.04,21.32,12.69,0.8,0.78,221.84,41.98,13.88,28.14,14.92,2.29,27.19,12.2,10.73,12.63,17.6,0.52,49.93,7.91,17.29),
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $rowVals = $data[$i]
  $targetRow = $i + 2
  for ($j = 0; $j -lt $rowVals.Length; $j++) {
    $ws.Cells.Item($targetRow, $j + 1).Value2 = $rowVals[$j]
  }
}

# Row 6 (the old 5th data row) is dropped entirely.
$ws.Rows.Item(6).Delete()

# A handful of columns got new autofit widths after the data refresh.
$narrowCols = @(10, 17, 27, 28, 30)
foreach ($c in $narrowCols) {
  $ws.Columns.Item($c).ColumnWidth = 6.166666666666667
}
$wideCols = @(20, 24)
foreach ($c in $wideCols) {
  $ws.Columns.Item($c).ColumnWidth = 7.166666666666667
}
